# Update column AF ("doctor_MA") stat values on Sheet1.
# Mapping of cell -> new value, per the commit's updated results.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "AF4"  = 0.755
    "AF5"  = 0.971
    "AF6"  = 0.849
    "AF7"  = 0.918
    "AF8"  = 0.88
    "AF9"  = 0.735
    "AF10" = 0.971
    "AF11" = 0.971
    "AF12" = 1.273
    "AF13" = 1.559
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
